# Apply updated cryptocurrency price/volume data per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.083.39"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.875.07"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("D5").Value = "'241.84"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.4886"
$ws.Range("E7").Value = "  -2.05%  "
$ws.Range("D8").Value = "'0.2899"
$ws.Range("E8").Value = "  -1.46%  "
$ws.Range("D9").Value = "'0.06575"
$ws.Range("E9").Value = "  -1.22%  "
$ws.Range("D10").Value = "1.878.20"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").Value = "'16.34"
$ws.Range("E11").Value = "  -3.86%  "
$ws.Range("D12").Value = "'0.07203"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "'0.6638"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").Value = "'4.903"
$ws.Range("E14").Value = "  +1.07%  "
$ws.Range("D15").Value = "'85.70"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("D16").Value = "30.056.87"
$ws.Range("E16").Value = "  -0.24%  "
$ws.Range("D17").Value = "'0.000007774"
$ws.Range("E17").Value = "  -2.99%  "
$ws.Range("D18").Value = "'0.9998"
$ws.Range("E18").Value = "  +0.33%  "
$ws.Range("D19").Value = "'12.71"
$ws.Range("E19").Value = "  -1.74%  "
$ws.Range("D20").Value = "2.119.70"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'4.754"
$ws.Range("E22").Value = "  -0.49%  "
$ws.Range("D23").Value = "'5.815"
$ws.Range("E23").Value = "  +2.71%  "
$ws.Range("D24").Value = "'9.161"
$ws.Range("E24").Value = "  -0.11%  "
$ws.Range("D25").Value = "'152.85"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("D26").Value = "'142.70"
$ws.Range("E26").Value = "  +7.02%  "
$ws.Range("D27").Value = "'16.93"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").Value = "'1.878"
$ws.Range("E28").Value = "  -3.64%  "
$ws.Range("E29").Value = "  +1.96%  "
$ws.Range("D30").Value = "'4.195"
$ws.Range("E30").Value = "  -0.80%  "
$ws.Range("D31").Value = "'0.08764"
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("D32").Value = "'3.995"
$ws.Range("E32").Value = "  +1.05%  "
$ws.Range("D33").Value = "'0.05119"
$ws.Range("E33").Value = "  -0.35%  "
$ws.Range("D34").Value = "'0.7129"
$ws.Range("E34").Value = "  +0.62%  "
$ws.Range("D35").Value = "'1.103"
$ws.Range("D36").Value = "'2.667"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("D37").Value = "'0.01827"
$ws.Range("E37").Value = "  +9.78%  "
$ws.Range("D38").Value = "'2.673"
$ws.Range("E38").Value = "  -3.66%  "
$ws.Range("D39").Value = "'2.126"
$ws.Range("E39").Value = "  -4.78%  "
$ws.Range("D40").Value = "'0.9231"
$ws.Range("E40").Value = "  -1.96%  "
$ws.Range("D41").Value = "'5.774"
$ws.Range("E41").Value = "  -4.94%  "
$ws.Range("D42").Value = "'0.9992"
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "'103.70"
$ws.Range("E43").Value = "  +0.56%  "
$ws.Range("D44").Value = "'0.4208"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "'7.386"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "'0.1278"
$ws.Range("E46").Value = "  +0.93%  "
$ws.Range("D47").Value = "'0.05708"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("D48").Value = "'32.76"
$ws.Range("E48").Value = "  -0.16%  "
$ws.Range("D49").Value = "'8.244"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").Value = "'0.3742"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("D51").Value = "'55.65"
$ws.Range("E51").Value = "  -0.60%  "
